$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to columns C (nombre_aides) and E (montant_total) for the
# 2022-06-01 refresh of the Fonds de solidarite / volet 1 regional data.
$updates = @(
    @{ Row = 2;   C = 766326;  E = 1429219098 },
    @{ Row = 3;   C = 791;     E = 2233378 },
    @{ Row = 48;  C = 150634;  E = 275739556 },
    @{ Row = 112; C = 145230;  E = 716365950 },
    @{ Row = 121; C = 1306266; E = 2275067309 },
    @{ Row = 129; C = 633607;  E = 3431371573 },
    @{ Row = 131; C = 378;     E = 19427930 },
    @{ Row = 132; C = 585849;  E = 3467997594 },
    @{ Row = 136; C = 26694;   E = 144319346 },
    @{ Row = 137; C = 51;      E = 2267833 },
    @{ Row = 154; C = 18462;   E = 73574082 },
    @{ Row = 186; C = 236828;  E = 1189788011 },
    @{ Row = 221; C = 135499;  E = 681875940 },
    @{ Row = 240; C = 205914;  E = 1069370264 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
